# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp text (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 23 de Octubre de 2020 a las 03:06"

# --- Swap "Montserrat" (row 216) and "Islas Malvinas" (row 217) entries,
#     including their per-row stats, so the two countries trade places
#     (mirrors the shared-string reorder in the source diff) ---
$a216 = $ws.Range("A216").Value()
$b216 = $ws.Range("B216").Value()
$c216 = $ws.Range("C216").Value()
$d216 = $ws.Range("D216").Value()
$e216 = $ws.Range("E216").Value()
$f216 = $ws.Range("F216").Value()
$g216 = $ws.Range("G216").Value()
$h216 = $ws.Range("H216").Value()

$a217 = $ws.Range("A217").Value()
$b217 = $ws.Range("B217").Value()
$c217 = $ws.Range("C217").Value()
$d217 = $ws.Range("D217").Value()
$e217 = $ws.Range("E217").Value()
$f217 = $ws.Range("F217").Value()
$g217 = $ws.Range("G217").Value()
$h217 = $ws.Range("H217").Value()

$ws.Range("A216").Value = $a217
$ws.Range("B216").Value = $b217
$ws.Range("C216").Value = $c217
$ws.Range("D216").Value = $d217
$ws.Range("E216").Value = $e217
$ws.Range("F216").Value = $f217
$ws.Range("G216").Value = $g217
$ws.Range("H216").Value = $h217

$ws.Range("A217").Value = $a216
$ws.Range("B217").Value = $b216
$ws.Range("C217").Value = $c216
$ws.Range("D217").Value = $d216
$ws.Range("E217").Value = $e216
$ws.Range("F217").Value = $f216
$ws.Range("G217").Value = $g216
$ws.Range("H217").Value = $h216

# --- Updated country statistics ---

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 8660399
$ws.Range("C4").Value = 73049
$ws.Range("D4").Value = 5654267
$ws.Range("E4").Value = 2777765
$ws.Range("G4").Value = 959
$ws.Range("H4").Value = 228367

# Row 66 - Paraguay
$ws.Range("B66").Value = 57526
$ws.Range("C66").Value = 707
$ws.Range("D66").Value = 38187
$ws.Range("E66").Value = 18077
$ws.Range("G66").Value = 12
$ws.Range("H66").Value = 1262

# Row 171 - San Marino
$ws.Range("B171").Value = 802
$ws.Range("C171").Value = 28
$ws.Range("D171").Value = 711
$ws.Range("E171").Value = 49

# Row 184 - Mauricio
$ws.Range("B184").Value = 425
$ws.Range("C184").Value = 6
$ws.Range("D184").Value = 386
$ws.Range("E184").Value = 29

# Row 193
$ws.Range("D193").Value = 175
$ws.Range("E193").Value = 4

# Row 197 - Bermudas
$ws.Range("D197").Value = 107
$ws.Range("E197").Value = 12

# Row 219 - Islas Salomon
$ws.Range("B219").Value = 4
$ws.Range("C219").Value = 1
$ws.Range("D219").Value = 3
$ws.Range("E219").Value = 1
